# Update "Backlog Taches01" sheet (Feuil1) with Scrum-poker results:
# Column C ("Efforts") now holds story-point-style effort values (can be
# fractional), and column D ("Minutes") is now a formula C*60 instead of a
# hard-coded literal.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New effort values for rows 3..21 (column C)
$efforts = @{
    3  = 1
    4  = 1
    5  = 0.5
    6  = 1
    7  = 0.5
    8  = 0.5
    9  = 2
    10 = 2
    11 = 1
    12 = 2
    13 = 0.5
    14 = 0.5
    15 = 0.5
    16 = 0.5
    17 = 0.5
    18 = 0.5
    19 = 1
    20 = 2
    21 = 3
}

foreach ($row in 3..21) {
    $ws.Cells.Item($row, 3).Value = $efforts[$row]
    $ws.Cells.Item($row, 3).Style = "Normal"
}

# Column D becomes a computed formula (minutes = effort hours * 60) instead
# of a hard-coded literal. Assign as two range-formulas (row 3 standalone,
# rows 4-21 as one relative-fill) so Excel builds the same shared-formula
# group the author ended up with.
$ws.Range("D3").Formula = "=C3*60"
$ws.Range("D4:D21").Formula = "=C4*60"

# View tweak recorded in the diff (selection moved while reviewing)
$ws.Range("O13").Select()
